# Add a new column "Vị trí kho" (warehouse location) right after the
# "ĐVT" column (column C) and before "Nhóm hàng" (old column D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D - this shifts old D..J to E..K and carries
# formatting from the surrounding cells, matching how Excel behaves when a
# column is inserted from the UI.
$ws.Columns.Item("D").Insert()

# Header text + style for the new header cell.
$ws.Range("D3").Value = "Vị trí kho"
$ws.Range("D3").Style = "Good"

# Match the header-row fill (light blue) style used across the whole header
# row after the edit - apply the themed fill to the header row band.
$headerRange = $ws.Range("A3:K3")
$headerRange.Interior.ThemeColor = 5
$headerRange.Interior.TintAndShade = 0.59999389629810485

# New column should mirror the width/format of its neighboring "ĐVT" column.
$ws.Columns.Item("D").ColumnWidth = $ws.Columns.Item("C").ColumnWidth

$wb.Save()
